$d = $word.ActiveDocument

# --- Part 1: merge the two "Fri Sep 14" / " 13:08:05 PDT 2017" runs into one run ---
# The visible text doesn't change, only the run split; Find/Replace across the
# run boundary collapses it into a single run, matching the target XML.
$d.Content.Find.Execute("Fri Sep 14 13:08:05 PDT 2017", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Fri Sep 14 13:08:05 PDT 2017", 2) | Out-Null

# --- Part 2: append the new "Sat Sep 15" purchase-details block ---
# Locate the paragraph that ends the "Fri Sep 14" block ("Amount balance ... - 313441.0")
# and insert the new paragraphs right after it, before the pre-existing trailing
# empty paragraphs.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*- 313441.0*") {
        $target = $cand
    }
}

$insertAt = $target.Range
$insertAt.Collapse(0)

$sep = "---------------------------------------------------------------"

$lines = @(
    "",
    "Sat Sep 15 12:52:52 PDT 2017",
    "Person Name`t`t`t`t- NSA",
    $sep,
    "Item Name`t`t`t`t- CARROT",
    "Number of Pockets`t`t`t- 4",
    "Number of KGs`t`t`t- 381",
    "Rate`t`t`t`t`t- 18",
    "Transport & Miscellaneous`t- 660",
    "Total Price`t`t`t`t- 7518.0",
    "Amount balance`t`t`t- 320959.0",
    "",
    ""
)

$block = ($lines -join "`r") + "`r"
$insertAt.InsertAfter($block)

# Re-fetch the paragraph index of $target (Paragraphs collection re-indexes live,
# $target itself still points at the same paragraph since text/content didn't move).
$startIndex = $target.Range.Information(3)  # wdActiveEndAdjustedPageNumber placeholder (unused)

$baseIndex = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $target.Range.Start) {
        $baseIndex = $i
        break
    }
}

# New paragraphs are baseIndex+1 .. baseIndex+13
$boldOffsets = @(1, 11, 13)
foreach ($off in $boldOffsets) {
    $d.Paragraphs.Item($baseIndex + $off).Range.Bold = 1
}

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
